# Re-sort the comma-separated "Recorded By" values in column G.
# Rule observed from the target diff: split the cell text on ", ", keep any
# literal lowercase "system" token pinned at the front (if present), and
# sort the remaining tokens using ordinal (case-sensitive) ascending order.
# This naturally places "System" (capital S) before lowercase email
# addresses, and sorts multiple email addresses alphabetically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ", "

    $pinned = @()
    $rest = @()
    foreach ($p in $parts) {
        if ($p -ceq "system") {
            $pinned += $p
        } else {
            $rest += $p
        }
    }

    $restSorted = $rest | Sort-Object { $_ } -Culture "en-US" -CaseSensitive

    $newParts = $pinned + $restSorted
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -cne $val) {
        $cell.Value = $newVal
    }
}
